# ---------------------------------------------------------------------------
# Commit: "Specific Baskets B6 and B7 for semester 7"
#   1. Re-assign the timetable slots in Section_A / Section_B to reflect the
#      newly scheduled baskets B6/B7, tagging every non-"Free" / non-break
#      slot with the classroom that was allocated to it ("[C405]").
#   2. Add a "Semester_Rules" sheet documenting the scheduling rule that was
#      applied.
#   3. Add a "Classroom_Utilization" sheet summarising room usage, showing
#      C405 now fully booked (56 weekly hours / 100% utilization).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1a. Section_A timetable updates
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$sectionA = @{
    "B2" = "ELECTIVE_B1 [C405]"; "C2" = "Free";                      "D2" = "ELECTIVE_B1 [C405]"; "E2" = "DS161 [C405]";             "F2" = "DS161 [C405]"
    "B3" = "EC161 [C405]";       "C3" = "Free";                      "D3" = "Free";                "E3" = "EC161 [C405]";             "F3" = "MA162 [C405]"
    "B5" = "Free";               "C5" = "Free";                      "D5" = "Free";                "E5" = "MA162 [C405]";             "F5" = "HS152 [C405]"
    "B6" = "CS161 (Tutorial) [C405]";  "C6" = "MA161 (Tutorial) [C405]"; "D6" = "HS152 (Tutorial) [C405]"; "E6" = "Free";              "F6" = "ELECTIVE_B1 (Tutorial) [C405]"
    "B7" = "MA161 [C405]";       "C7" = "HS152 [C405]";              "D7" = "CS161 [C405]";        "E7" = "MA161 [C405]";             "F7" = "CS161 [C405]"
    "B8" = "DS161 (Tutorial) [C405]";  "C8" = "MA162 (Tutorial) [C405]"; "D8" = "EC161 (Tutorial) [C405]"; "E8" = "Free";              "F8" = "Free"
}
foreach ($addr in $sectionA.Keys) {
    $wsA.Range($addr).Value = $sectionA[$addr]
}

# ---------------------------------------------------------------------------
# 1b. Section_B timetable updates
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$sectionB = @{
    "B2" = "ELECTIVE_B1 [C405]"; "C2" = "HS152 [C405]";              "D2" = "ELECTIVE_B1 [C405]"; "E2" = "MA162 [C405]";             "F2" = "DS161 [C405]"
    "B3" = "MA162 [C405]";       "C3" = "Free";                      "D3" = "DS161 [C405]";       "E3" = "Free";                      "F3" = "CS161 [C405]"
    "B5" = "Free";               "C5" = "Free";                      "D5" = "EC161 [C405]";       "E5" = "Free";                      "F5" = "MA161 [C405]"
    "B6" = "Free";               "C6" = "MA161 (Tutorial) [C405]";   "D6" = "CS161 (Tutorial) [C405]"; "E6" = "EC161 (Tutorial) [C405]"; "F6" = "ELECTIVE_B1 (Tutorial) [C405]"
    "B7" = "CS161 [C405]";       "C7" = "Free";                      "D7" = "MA161 [C405]";       "E7" = "HS152 [C405]";              "F7" = "EC161 [C405]"
    "B8" = "Free";               "C8" = "DS161 (Tutorial) [C405]";   "D8" = "MA162 (Tutorial) [C405]"; "E8" = "Free";                  "F8" = "HS152 (Tutorial) [C405]"
}
foreach ($addr in $sectionB.Keys) {
    $wsB.Range($addr).Value = $sectionB[$addr]
}

# ---------------------------------------------------------------------------
# 2. New sheet: Semester_Rules
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRules = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsRules.Name = "Semester_Rules"

$wsRules.Range("A1").Value = "Semester"
$wsRules.Range("B1").Value = "Rule"
$wsRules.Range("C1").Value = "Exclusion"
$wsRules.Range("D1").Value = "Reason"
$wsRules.Range("E1").Value = "Scheduled Baskets"
$wsRules.Range("F1").Value = "Status"

$wsRules.Range("A2").Value = "Semester 1"
$wsRules.Range("B2").Value = "Schedule all elective baskets"
$wsRules.Range("C2").Value = "None"
$wsRules.Range("D2").Value = "No specific restrictions for this semester"
$wsRules.Range("E2").Value = "ELECTIVE_B1"
$wsRules.Range("F2").Value = [char]0x2705 + " Applied"

# Re-use the bold/centered/bordered header style already present in the
# workbook (e.g. Course_Summary!A1) for the new header row.
$wsCourseSummary = $wb.Worksheets.Item("Course_Summary")
$wsCourseSummary.Range("A1").Copy()
$wsRules.Range("A1:F1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. New sheet: Classroom_Utilization
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRooms = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsRooms.Name = "Classroom_Utilization"

$wsRooms.Range("A1").Value = "Room Number"
$wsRooms.Range("B1").Value = "Type"
$wsRooms.Range("C1").Value = "Capacity"
$wsRooms.Range("D1").Value = "Weekly Hours (Timetable)"
$wsRooms.Range("E1").Value = "Daily Avg Hours (Timetable)"
$wsRooms.Range("F1").Value = "Exam Sessions"
$wsRooms.Range("G1").Value = "Utilization Rate (%)"
$wsRooms.Range("H1").Value = "Facilities"

$wsCourseSummary.Range("A1").Copy()
$wsRooms.Range("A1:H1").PasteSpecial(-4122)

$rooms = @(
    @("C001","Recreation","nil",0,0,0,0,""),
    @("C002","large classroom","116",0,0,0,0,"Projector"),
    @("C003","large classroom","135",0,0,0,0,"Projector"),
    @("C004","Auditorium","240",0,0,0,0,"Projector"),
    @("C101","classroom","96",0,0,0,0,"Projector"),
    @("C102","classroom","96",0,0,0,0,"Projector"),
    @("C103","library","nil",0,0,0,0,"Computers"),
    @("C104","classroom","96",0,0,0,0,"Projector"),
    @("L105","Hardware Lab","40",0,0,0,0,"Hardware Equipment"),
    @("L106","Software Lab","40",0,0,0,0,"Computers"),
    @("L107","Software Lab","40",0,0,0,0,"Computers"),
    @("C201","classroom","96",0,0,0,0,"Projector"),
    @("C202","classroom","96",0,0,0,0,"Projector"),
    @("C203","classroom","96",0,0,0,0,"Projector"),
    @("C204","classroom","96",0,0,0,0,"Projector"),
    @("C205","classroom","96",0,0,0,0,"Projector"),
    @("L206","Hardware Lab","40",0,0,0,0,"Hardware Equipment"),
    @("L207","Software Lab","40",0,0,0,0,"Computers"),
    @("L208","Software Lab","40",0,0,0,0,"Computers"),
    @("C301","Physics Lab","40",0,0,0,0,"Projector"),
    @("C302","classroom","96",0,0,0,0,"Projector"),
    @("C303","classroom","96",0,0,0,0,"Projector"),
    @("C304","classroom","96",0,0,0,0,"Projector"),
    @("C305","classroom","96",0,0,0,0,"Projector"),
    @("L306","classroom","96",0,0,0,0,"Computers"),
    @("L307","Research Scholar Lab","40",0,0,0,0,"Computers"),
    @("L308","Research Scholar Lab","40",0,0,0,0,"Computers"),
    @("C401","classroom","96",0,0,0,0,"Projector"),
    @("C402","classroom","96",0,0,0,0,"Projector"),
    @("C403","classroom","78",0,0,0,0,"Projector"),
    @("C404","classroom","78",0,0,0,0,"Projector"),
    @("C405","classroom","78",56,11.2,0,100,"Projector"),
    @("L406","classroom","78",0,0,0,0,"Computers"),
    @("L407","classroom","78",0,0,0,0,"Computers"),
    @("L408","classroom without projector","78",0,0,0,0,"Computers")
)

# The "Capacity" column holds numeric-looking text (e.g. "96", "nil") in the
# source data, so force text formatting on that column before writing so the
# numeric-looking values aren't silently coerced into numbers.
$wsRooms.Range("C3:C36").NumberFormat = "@"

$r = 2
foreach ($room in $rooms) {
    $wsRooms.Cells.Item($r, 1).Value = $room[0]
    $wsRooms.Cells.Item($r, 2).Value = $room[1]
    $wsRooms.Cells.Item($r, 3).Value = $room[2]
    $wsRooms.Cells.Item($r, 4).Value = $room[3]
    $wsRooms.Cells.Item($r, 5).Value = $room[4]
    $wsRooms.Cells.Item($r, 6).Value = $room[5]
    $wsRooms.Cells.Item($r, 7).Value = $room[6]
    $wsRooms.Cells.Item($r, 8).Value = $room[7]
    $r++
}

# ---------------------------------------------------------------------------
# Leave the workbook focused on the first sheet, as it was originally.
# ---------------------------------------------------------------------------
$wsA.Activate()
